# Update cryptocurrency price/volume data to latest scraped values
# (mirrors the GitHub Actions "Updated cryptos list" commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.709.09"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.894.90"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  -1.29%  "
$ws.Range("D5").Value = "'312.58"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("D7").Value = "'0.4921"
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("D8").Value = "'0.3798"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").Value = "'0.07328"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").Value = "'0.9142"
$ws.Range("E10").Value = "  -2.74%  "
$ws.Range("D11").Value = "'20.58"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07673"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.897.62"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "'5.475"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "'6.603"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "'91.03"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").Value = "27.931.77"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "'14.53"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").Value = "'5.128"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "2.169.61"
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").Value = "'1.912"
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("D26").Value = "'153.66"
$ws.Range("E26").Value = "  -2.15%  "
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("D28").Value = "'2.147"
$ws.Range("E28").Value = "  +4.71%  "
$ws.Range("D29").Value = "'115.71"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "'4.899"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").Value = "'0.08933"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").Value = "'3.197"
$ws.Range("E32").Value = "  -4.15%  "
$ws.Range("D33").Value = "'1.221"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").Value = "'0.7672"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'4.638"
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").Value = "'2.530"
$ws.Range("E37").Value = "  -7.64%  "
$ws.Range("E38").Value = "  -3.57%  "
$ws.Range("D39").Value = "'0.05270"
$ws.Range("D40").Value = "'0.5462"
$ws.Range("E40").Value = "  -2.50%  "
$ws.Range("D41").Value = "'2.979"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("D42").Value = "'6.911"
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("D43").Value = "'8.529"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").Value = "'112.45"
$ws.Range("E44").Value = "  +6.81%  "
$ws.Range("D45").Value = "'0.1518"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").Value = "'10.59"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "'0.4790"
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").Value = "'1.631"
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D50").Value = "'67.41"
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("D51").Value = "'0.06048"
$ws.Range("E51").Value = "  -1.17%  "
